$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-10 (columns A:T), reflecting the updated
# NATMI Il16-Cd4 LR-pair statistics (re-run "following Dr Hou advice").
$data = @(
    @("ECs", "Il16", "Cd4", "ECs", 3, 1, 8.267595333333334, 24.802786, 0.3656664502891759, 0.3656664502891758, 2, 0.6666666666666666, 0.539832, 1.619496, 0.3224687287194008, 0.3224687287194007, 4.463112523984, 40.168012715856, 0.1179159953600865, 0.1179159953600865),
    @("ECs", "Il16", "Cd4", "FAPs", 3, 1, 8.267595333333334, 24.802786, 0.3656664502891759, 0.3656664502891758, 3, 1, 1.034252666666666, 3.102758, 0.6178109904463799, 0.6178109904463799, 8.550782520420888, 76.95704268378799, 0.2259127518261677, 0.2259127518261677),
    @("ECs", "Il16", "Cd4", "sCs", 3, 1, 8.267595333333334, 24.802786, 0.3656664502891759, 0.3656664502891758, 2, 0.6666666666666666, 0.09997533333333335, 0.299926, 0.05972028083421942, 0.05972028083421941, 0.8265555993151114, 7.439000393836001, 0.02183770310292172, 0.02183770310292171),
    @("FAPs", "Il16", "Cd4", "ECs", 3, 1, 11.10830433333333, 33.324913, 0.491307816908375, 0.491307816908375, 2, 0.6666666666666666, 0.539832, 1.619496, 0.3224687287194008, 0.3224687287194007, 5.996618144871999, 53.96956330384799, 0.1584314071283478, 0.1584314071283478),
    @("FAPs", "Il16", "Cd4", "FAPs", 3, 1, 11.10830433333333, 33.324913, 0.491307816908375, 0.491307816908375, 3, 1, 1.034252666666666, 3.102758, 0.6178109904463799, 0.6178109904463799, 11.48879337889488, 103.399140410054, 0.3035353689782118, 0.3035353689782118),
    @("FAPs", "Il16", "Cd4", "sCs", 3, 1, 11.10830433333333, 33.324913, 0.491307816908375, 0.491307816908375, 2, 0.6666666666666666, 0.09997533333333335, 0.299926, 0.05972028083421942, 0.05972028083421941, 1.110556428493111, 9.995007856437999, 0.02934104080181541, 0.0293410408018154),
    @("sCs", "Il16", "Cd4", "ECs", 3, 1, 3.233763666666666, 9.701290999999999, 0.1430257328024492, 0.1430257328024492, 2, 0.6666666666666666, 0.539832, 1.619496, 0.3224687287194008, 0.3224687287194007, 1.745689107704, 15.711201969336, 0.04612132623096651, 0.0461213262309665),
    @("sCs", "Il16", "Cd4", "FAPs", 3, 1, 3.233763666666666, 9.701290999999999, 0.1430257328024492, 0.1430257328024492, 3, 1, 1.034252666666666, 3.102758, 0.6178109904463799, 0.6178109904463799, 3.344528695619777, 30.10075826057799, 0.08836286964200045, 0.08836286964200045),
    @("sCs", "Il16", "Cd4", "ECs", 3, 1, 3.233763666666666, 9.701290999999999, 0.1430257328024492, 0.1430257328024492, 2, 0.6666666666666666, 0.09997533333333335, 0.299926, 0.05972028083421942, 0.05972028083421941, 0.3232966004962222, 2.909669404466, 0.008541536929482297, 0.008541536929482296)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}
